# Apply odds updates to the "Jogos da Semana" workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 9
$ws.Range("O2").Value = 1.36
$ws.Range("P2").Value = 3
$ws.Range("Q2").Value = 2.15
$ws.Range("R2").Value = 1.67

# --- Row 4 updates ---
$ws.Range("J4").Value = 3.2
$ws.Range("L4").Value = 3.25
$ws.Range("W4").Value = 8.5
$ws.Range("X4").Value = 13.5
$ws.Range("Y4").Value = 9.5
$ws.Range("AA4").Value = 21
$ws.Range("AB4").Value = 29
$ws.Range("AC4").Value = 9.5
$ws.Range("AH4").Value = 8.75
$ws.Range("AI4").Value = 14
$ws.Range("AL4").Value = 22
$ws.Range("AM4").Value = 29
$ws.Range("AN4").Value = 4.5
$ws.Range("AO4").Value = 14
$ws.Range("AP4").Value = 21
$ws.Range("AQ4").Value = 65
$ws.Range("AR4").Value = 100
$ws.Range("AS4").Value = 250
$ws.Range("AT4").Value = 2.45
$ws.Range("AX4").Value = 14.5
$ws.Range("AZ4").Value = 65
$ws.Range("BA4").Value = 100
$ws.Range("BB4").Value = 250
